$d = $word.ActiveDocument

# Paragraph 1
$p = $d.Paragraphs(1)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">This grammar is lightly adapted from the provided grammar. I refactored both declaration-list and param to convert this from a LL3 grammar to a LL1 grammar. (Before, the rule was param := type-specifier ID [] | type-specifier ID ε. By refactoring, we no longer have </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>too</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> look ahead three items to tokens determine which rule applies here).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 6
$p = $d.Paragraphs(6)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>type-specifier := int | void</w:t></w:r><w:r><w:t xml:space="preserve"> | float</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 7
$p = $d.Paragraphs(7)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>fun-declaration := ( params ) compound-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 11
$p = $d.Paragraphs(11)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>compound-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> := { local-declarations statement-list }</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 14
$p = $d.Paragraphs(14)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>statement := expression-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> | compound-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> | selection-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> | iteration-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> | return-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> | io-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 15
$p = $d.Paragraphs(15)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>io-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> := input-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> | output-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 16
$p = $d.Paragraphs(16)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>input-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> := input ( STRING )</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 17
$p = $d.Paragraphs(17)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>output-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> := output ( STRING ) | output ( expression )</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 18
$p = $d.Paragraphs(18)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>expression-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> := expression ; | ;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 19
$p = $d.Paragraphs(19)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>selection-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> := if ( simple-expression ) statement | if ( simple-expression ) statement else statement</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 20
$p = $d.Paragraphs(20)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>iteration-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> := while ( expression ) statement</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 21
$p = $d.Paragraphs(21)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>return-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> := return ; | return expression ;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 24
$p = $d.Paragraphs(24)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">simple-expression := additive-expression </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>relop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> additive-expression | additive-expression</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 25
$p = $d.Paragraphs(25)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>relop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> := &lt;= | &lt; | &gt; | &gt;= | == | !=</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 26
$p = $d.Paragraphs(26)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">additive-expression := additive-expression </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>addop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> term | term</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 27
$p = $d.Paragraphs(27)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>addop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> := + | -</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 28
$p = $d.Paragraphs(28)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">term := term </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mulop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> factor | factor</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 29
$p = $d.Paragraphs(29)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>mulop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> := * | /</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 30
$p = $d.Paragraphs(30)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">factor := ( simple-expression ) | var | call | NUM | </w:t></w:r><w:r><w:t xml:space="preserve">FLOAT | </w:t></w:r><w:r><w:t>input-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stmt</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 31
$p = $d.Paragraphs(31)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">call := ID ( </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>args</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> )</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 32
$p = $d.Paragraphs(32)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>args</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> := </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>arg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-list | ε</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# Paragraph 33
$p = $d.Paragraphs(33)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>arg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-list := </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>arg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-list , expression | expression</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

Write-Host "done"
